$d = $word.ActiveDocument

# 1) "Actores" table cell: "Usuario registrado, Usuario no registrado" -> "Usuario registrado"
$d.Content.Find.Execute("Usuario registrado, Usuario no registrado", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Usuario registrado", 2)

# 2) Remove spurious spell-check proofErr markers around "Nº" (no visible text change,
#    but the proofErr elements should be removed). Re-setting the text of the run
#    forces the proofErr wrapper to be dropped.
$d.Content.Find.Execute("Nº", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Nº", 2)
